$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Coord: normal vector scan"
$ws.Range("C1").Style = $ws.Range("B1").Style

# Updated B values (slightly refined precision) and new C values (vector strings)
$data = @(
    @{Row=2;  B=0.1012942899029146;  C="[0.         0.32329375 0.94629866]"},
    @{Row=3;  B=1.056430168733275;   C="[-0.43976721  0.52808933  0.72644784]"},
    @{Row=4;  B=0.3648215635934389;  C="[-0.00491449  0.02008922  0.99978611]"},
    @{Row=5;  B=0.7774609747866914;  C="[-0.00146906  0.28154281 -0.95954754]"},
    @{Row=6;  B=2.049911883240604;   C="[0.72088551 0.28601219 0.63128528]"},
    @{Row=7;  B=0.9772765525014188;  C="[-0.73289099 -0.26794639  0.62536031]"},
    @{Row=8;  B=0.7762933880188216;  C="[0.         0.31212327 0.95004161]"},
    @{Row=9;  B=0.9539078768857674;  C="[ 0.         -0.3091767   0.95100461]"},
    @{Row=10; B=1.296331392875424;   C="[-0.73439777  0.25940504  0.62719131]"},
    @{Row=11; B=2.032655723420702;   C="[ 0.72554665 -0.27723529  0.62985923]"},
    @{Row=12; B=0.7423662776918968;  C="[ 0.00144528 -0.28213156 -0.95937464]"},
    @{Row=13; B=2.452203874206384;   C="[-0.718175   -0.28106773  0.63657332]"},
    @{Row=14; B=4.120309090252077;   C="[0.69963041 0.25511099 0.66740968]"},
    @{Row=15; B=1.294056002612753;   C="[ 0.         -0.30352545  0.95282333]"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
